$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.420.82'
$ws.Range("E2").Value = '  +1.79%  '

$ws.Range("D3").Value = '2.163.03'
$ws.Range("E3").Value = '  +2.97%  '

$ws.Range("E4").Value = '  -0.07%  '

$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '227.90'
$r.Style = "Normal"
$ws.Range("E5").Value = '  -0.23%  '

$ws.Range("E6").Value = '  +1.17%  '

$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = '63.81'
$r.Style = "Normal"
$ws.Range("E7").Value = '  +2.97%  '

$ws.Range("E9").Value = '  +1.75%  '

$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = '0.0854'
$r.Style = "Normal"
$ws.Range("E10").Value = '  +1.16%  '

$ws.Range("E11").Value = '  -0.10%  '

$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = '16.05'
$r.Style = "Normal"
$ws.Range("E12").Value = '  +1.47%  '

$ws.Range("D13").Value = '2.482.34'
$ws.Range("E13").Value = '  +2.96%  '

$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = '22.13'
$r.Style = "Normal"
$ws.Range("E14").Value = '  +0.31%  '

$ws.Range("E15").Value = '  +1.00%  '

$ws.Range("E16").Value = '  +0.23%  '

$ws.Range("D17").Value = '2.156.22'
$ws.Range("E17").Value = '  +3.07%  '

$ws.Range("D18").Value = '39.394.56'
$ws.Range("E18").Value = '  +1.71%  '

$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = '6.16'
$r.Style = "Normal"
$ws.Range("E19").Value = '  +1.19%  '

$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = '71.89'
$r.Style = "Normal"
$ws.Range("E20").Value = '  +0.08%  '

$ws.Range("E21").Value = '  +1.15%  '

$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = '229.70'
$r.Style = "Normal"
$ws.Range("E22").Value = '  +0.97%  '

$ws.Range("E23").Value = '  +0.05%  '

$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = '2.35'
$r.Style = "Normal"
$ws.Range("E24").Value = '  -1.10%  '

$ws.Range("E25").Value = '  +1.76%  '

$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = '9.69'
$r.Style = "Normal"
$ws.Range("E26").Value = '  +1.32%  '

$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = '172.18'
$r.Style = "Normal"
$ws.Range("E27").Value = '  +0.07%  '

$ws.Range("E28").Value = '  +2.17%  '

$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = '19.88'
$r.Style = "Normal"
$ws.Range("E29").Value = '  +2.75%  '

$ws.Range("E30").Value = '  -0.05%  '

$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = '2.61'
$r.Style = "Normal"
$ws.Range("E31").Value = '  +2.39%  '

$ws.Range("E32").Value = '  +1.23%  '

$ws.Range("E33").Value = '  +1.60%  '

$ws.Range("E34").Value = '  +1.85%  '

$ws.Range("E35").Value = '  -0.81%  '

$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = '0.0619'
$r.Style = "Normal"
$ws.Range("E36").Value = '  +0.09%  '

$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = '2.44'
$r.Style = "Normal"
$ws.Range("E37").Value = '  +0.63%  '

$ws.Range("E38").Value = '  -0.79%  '

$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = '0.999'
$r.Style = "Normal"
$ws.Range("E39").Value = '  -0.18%  '

$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = '103.28'
$r.Style = "Normal"
$ws.Range("E40").Value = '  +0.50%  '

$ws.Range("E41").Value = '  +0.46%  '

$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = '17.88'
$r.Style = "Normal"
$ws.Range("E42").Value = '  -1.26%  '

$ws.Range("D43").Value = '1.522.60'
$ws.Range("E43").Value = '  -0.82%  '

$ws.Range("E44").Value = '  +3.94%  '

$ws.Range("B45").Value = 'ARBITRUM'
$ws.Range("C45").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = '1.11'
$r.Style = "Normal"
$ws.Range("E45").Value = '  +5.54%  '

$ws.Range("E46").Value = '  +0.67%  '

$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = '0.0926'
$r.Style = "Normal"
$ws.Range("E47").Value = '  +1.67%  '

$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = '4.27'
$r.Style = "Normal"
$ws.Range("E48").Value = '  +3.58%  '

$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = '7.73'
$r.Style = "Normal"
$ws.Range("E49").Value = '  -1.33%  '

$ws.Range("D50").Value = '2.364.92'
$ws.Range("E50").Value = '  +3.11%  '

$ws.Range("E51").Value = '  -0.42%  '
